$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column C (copy style from B1 header, which is bold w/ border)
$ws.Range("C1").Value = "min_units"
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)

# min_units values for rows 2..68 (course rows)
$values = @(12,12,12,12,12,12,12,12,12,12,12,12,12,12,12,12,12,12,12,12,12,12,12,12,12,9,9,9,9,9,9,9,9,9,9,12,12,12,12,12,12,10,10,10,10,10,10,12,12,12,11,11,9,9,9,9,9,9,9,9,9,9,12,12,12,3,3)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $values[$i]
}
